# Weekly update: insert the latest week's two records (Primera/Segunda)
# at the top of the data block and push the rest of the history down.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at row 8, pushing existing rows 8-112 down to 10-114.
$ws.Range("A8:A9").EntireRow.Insert()

# New row 8: Caigua - Primera
$ws.Cells.Item(8, 1).Value = 1
$ws.Cells.Item(8, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(8, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(8, 4).Value = 44750
$ws.Cells.Item(8, 5).Value = 15
$ws.Cells.Item(8, 6).Value = 100112036
$ws.Cells.Item(8, 7).Value = "Caigua"
$ws.Cells.Item(8, 8).Value = "Sin especificar"
$ws.Cells.Item(8, 9).Value = "Primera"
$ws.Cells.Item(8, 10).Value = 130
$ws.Cells.Item(8, 11).Value = 8000
$ws.Cells.Item(8, 12).Value = 9000
$ws.Cells.Item(8, 13).Value = 8500
$ws.Cells.Item(8, 14).Value = "`$/caja 20 kilos"
$ws.Cells.Item(8, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(8, 16).Value = 425
$ws.Cells.Item(8, 17).Value = 20
$ws.Cells.Item(8, 18).Value = "Hortaliza"

# New row 9: Caigua - Segunda
$ws.Cells.Item(9, 1).Value = 1
$ws.Cells.Item(9, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(9, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(9, 4).Value = 44750
$ws.Cells.Item(9, 5).Value = 15
$ws.Cells.Item(9, 6).Value = 100112036
$ws.Cells.Item(9, 7).Value = "Caigua"
$ws.Cells.Item(9, 8).Value = "Sin especificar"
$ws.Cells.Item(9, 9).Value = "Segunda"
$ws.Cells.Item(9, 10).Value = 120
$ws.Cells.Item(9, 11).Value = 7000
$ws.Cells.Item(9, 12).Value = 8000
$ws.Cells.Item(9, 13).Value = 7500
$ws.Cells.Item(9, 14).Value = "`$/caja 20 kilos"
$ws.Cells.Item(9, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(9, 16).Value = 375
$ws.Cells.Item(9, 17).Value = 20
$ws.Cells.Item(9, 18).Value = "Hortaliza"

# Make sure the date column keeps the date number format used throughout
# column D (inherited from the row-insert, but set explicitly to be safe).
$ws.Range("D8:D9").NumberFormat = $ws.Range("D10").NumberFormat
